$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the little "Param Sims / True Value" side table (G1:H4) entirely.
$ws.Range("G1:H4").Clear()

# --- R_flat_prior block (rows 2-4): rename free-param "r1" -> "ma1" and
#     change its prior from Log-Normal(0, 5) to Uniform(0, 1).
$ws.Range("C3").Value = "ma1"
$ws.Range("D3").Value = "Uniform(0, 1)"

# --- Block at rows 6-8: R_intermedMa_prior -> R_strongMa_prior, with the
#     free param renamed r1/ma2 -> ma1 and the priors strengthened.
$ws.Range("A6").Value = "R_strongMa_prior"
$ws.Range("C7").Value = "ma1"
$ws.Range("D7").Value = "Beta(100, 900)"
$ws.Range("D8").Value = "Beta(500, 500)"

# --- Block at rows 10-12: R_intermedI0_prior -> R_strongI0_prior, with the
#     I0 prior strengthened and the free param renamed r1 -> ma1.
$ws.Range("A10").Value = "R_strongI0_prior"
$ws.Range("D10").Value = "Log-Normal(0.69, 0.05)"
$ws.Range("C11").Value = "ma1"
$ws.Range("D11").Value = "Uniform(0, 1)"

# --- Block at rows 14-18: R_strongMa_prior -> R_reparameterize_Ma_prior,
#     now re-parameterized with four free params (r1, r2, r3, ma4), each
#     Uniform(0, 1), and the I0 prior switched to Log-Normal(0.69, 0.05).
$ws.Range("A14").Value = "R_reparameterize_Ma_prior"
$ws.Range("D14").Value = "Log-Normal(0.69, 0.05)"
$ws.Range("C15").Value = "r1"
$ws.Range("D15").Value = "Uniform(0, 1)"
$ws.Range("C16").Value = "r2"
$ws.Range("D16").Value = "Uniform(0, 1)"
$ws.Range("C17").Value = "r3"
$ws.Range("D17").Value = "Uniform(0, 1)"
$ws.Range("C18").Value = "ma4"
$ws.Range("D18").Value = "Uniform(0, 1)"

# The old R_strongI0_prior block used to start its header at row 18
# (A18/B18) -- that header moved up to row 14, so clear the leftover A18/B18
# cells, then drop the leftover free-param rows 19-20 entirely.
$ws.Range("A18:B18").Clear()
$ws.Rows("19:20").Delete()

# Column A needs to be a bit wider for "R_reparameterize_Ma_prior".
# (ColumnWidth is stored with a +5/6 offset in the saved XML width, so back
# that off here to land exactly on width="27".)
$ws.Columns("A").ColumnWidth = 27 - (5/6)

# Leave the selection where editing ended up.
$ws.Range("D17").Select()
